# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the zh-cn and
# de-de report rows, bumping the handback run to a later point in time.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 12:43:35"
$wsZhCn.Range("H2").Value = "2016-03-11 12:43:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 12:43:39"
$wsDeDe.Range("H2").Value = "2016-03-11 12:43:57"
